# Security Program Tracking Template - SIEM sheet update
# Inserts a new "Sub-Items" column (B) on the SIEM sheet, expands the
# Description and Improvement Plan text for every capability row, and
# widens the affected columns accordingly.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("SIEM")

# ---------------------------------------------------------------------
# 1. Insert a new column B ("Sub-Items"); existing B..G shift to C..H
# ---------------------------------------------------------------------
$ws.Columns.Item(2).Insert()

# ---------------------------------------------------------------------
# 2. Header row
# ---------------------------------------------------------------------
$ws.Range("B1").Value = "Sub-Items"

# ---------------------------------------------------------------------
# 3. Sub-Items column values (new column B)
# ---------------------------------------------------------------------
$ws.Range("B2").Value = "Malware Detection, Anomaly Detection, Insider Threat Detection"
$ws.Range("B3").Value = "IOC Search, Behavioral Analysis, Threat Hypothesis Testing"
$ws.Range("B4").Value = "Executive Dashboard, Analyst Dashboard, Incident Response Dashboard"
$ws.Range("B5").Value = "Rule-Based Alerts, Threshold Alerts, Escalation Workflows, Machine Learning Alerts"
$ws.Range("B6").Value = "User Management, Log Retention, Role-Based Access Control, Configuration Auditing"
$ws.Range("B7").Value = "Software Updates, Performance Monitoring, Health Checks, Backup & Recovery"
$ws.Range("B8").Value = "Feed Integration, Threat Correlation, Automated Enrichment, Threat Actor Profiling"

# ---------------------------------------------------------------------
# 4. Expanded Description text (column C, formerly B)
# ---------------------------------------------------------------------
$ws.Range("C2").Value = "Ability to detect various threats, such as malware, anomalies, and insider threats using the SIEM. Includes correlation of data from multiple sources to identify complex attack patterns."
$ws.Range("C3").Value = "Ability to proactively search for indicators of compromise (IOCs) and unknown threats. Involves using both automated tools and manual analysis to detect stealthy or sophisticated attacks."
$ws.Range("C4").Value = "Customizable dashboards to monitor specific activities or use cases. Provides detailed views for different stakeholders, such as executives, analysts, and incident responders."
$ws.Range("C5").Value = "Ability to create and manage alerts for suspicious activities detected in logs. Includes advanced customization of alert logic, machine learning-based alerts, and escalation workflows based on severity."
$ws.Range("C6").Value = "Administrative tasks such as user management, log retention, configuration changes, and auditing configurations. Includes role-based access control, audit logging, and periodic review of permissions."
$ws.Range("C7").Value = "Regular maintenance, software updates, and system health checks for the SIEM. Includes monitoring for performance bottlenecks, ensuring timely software patches, and backup & recovery processes."
$ws.Range("C8").Value = "Capability to search for and correlate threat intelligence feeds with internal activity to identify threats. Supports integration with multiple threat feeds, automated enrichment of indicators, and threat actor profiling."

# ---------------------------------------------------------------------
# 5. Expanded Improvement Plan text (column G, formerly F)
# ---------------------------------------------------------------------
$ws.Range("G2").Value = "Improve detection rules for insider threats; Enhance correlation capabilities"
$ws.Range("G3").Value = "Increase frequency of threat hunts; Develop more advanced threat hunting playbooks"
$ws.Range("G4").Value = "Develop more user-specific dashboards; Create dashboards for executive reporting"
$ws.Range("G5").Value = "Refine alert thresholds to reduce false positives; Implement alert prioritization based on risk"
$ws.Range("G6").Value = "Streamline user role assignment; Automate audit logging for configuration changes"
$ws.Range("G7").Value = "Automate routine maintenance tasks; Implement proactive health checks; Ensure regular backup and recovery tests"
$ws.Range("G8").Value = "Integrate additional threat feeds; Automate enrichment of threat intelligence data; Expand threat actor profiling capabilities"

# ---------------------------------------------------------------------
# 6. Column widths - only touch the columns whose width actually changed
# ---------------------------------------------------------------------
$ws.Columns.Item(2).ColumnWidth = 69.83333333   # B - Sub-Items            -> 70.6640625 (closest achievable)
$ws.Columns.Item(3).ColumnWidth = 175.33333333  # C - Description          -> 176.1640625 (closest achievable)
$ws.Columns.Item(7).ColumnWidth = 101.66666667  # G - Improvement Plan     -> 102.5

# ---------------------------------------------------------------------
# 7. Select the full sheet (matches the saved selection state in the diff)
# ---------------------------------------------------------------------
$ws.Cells.Select()
